$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leaderboard")

# Column C holds the "Weighted Total" as text (e.g. "65.50"), so force
# text formatting before writing values to avoid Excel auto-converting
# them to numbers (which would drop the trailing zero).
$ws.Range("C2:C11").NumberFormat = "@"

# Update existing rows 2-5 with new values
$ws.Range("B2").Value = "Hacktronics"
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 8
$ws.Range("G2").Value = "Animal Type Classification - DEVANG SHUKLA.pdf"

$ws.Range("B3").Value = "Alt-Era"
$ws.Range("C3").Value = "64.00"
$ws.Range("D3").Value = 7
$ws.Range("G3").Value = "Alt-Era - KRISH PATHAK.pdf"

$ws.Range("B4").Value = "BenzeneCoder"
$ws.Range("C4").Value = "64.00"
$ws.Range("D4").Value = 7
$ws.Range("G4").Value = "BenzeneCoder - SARTHAK TIWARI.pdf"

$ws.Range("B5").Value = "Angaari Paltan"
$ws.Range("C5").Value = "63.50"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = "Angaari Paltan - ISHITA GOYAL.pdf"

# New row 6: Coding Pirates (previously row 3 content)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Coding Pirates"
$ws.Range("C6").Value = "61.50"
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = "2025CodingPirates - YASH KASAUDHAN.pdf"

# New row 7: Kairos
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Kairos"
$ws.Range("C7").Value = "60.00"
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 8
$ws.Range("G7").Value = "AgriNiti - TEENA gla.pdf"

# New row 8: Algo Wizards (previously row 5 content)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Algo Wizards"
$ws.Range("C8").Value = "58.50"
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 7
$ws.Range("G8").Value = "Algo wizards - LAXMI gla.pdf"

# New row 9: Binary Brains
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Binary Brains"
$ws.Range("C9").Value = "55.50"
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 7
$ws.Range("G9").Value = "Binary Brains - Milan Sharma.pdf"

# New row 10: ALT_F4
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "ALT_F4"
$ws.Range("C10").Value = "54.00"
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 7
$ws.Range("G10").Value = "ALT_f4 - VAIBHAV KUMAR.pdf"

# Row 11: AlgoYoddhas, moved from old row 6, rank now 10
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "AlgoYoddhas"
$ws.Range("C11").Value = "46.00"
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 6
$ws.Range("G11").Value = "AlgoYoddhas - ADWAIT PATEL.pdf"

# Clear the old row 6 leftover cell content beyond column G if any (not needed, row 6 now fully redefined above)
